# Add Betting Markets Analytics (Handicap, O/U, BTTS)
#
# The two stat rows for Napoli's centre-back pairing were re-ordered:
# row 19 ("Mathías Olivera") and row 20 ("Amir Rrahmani") swap all of
# their per-player statistics (columns C:DJ), while League/Team
# (A:B) and type/goalsPrevented (DK:DL) stay put since they are
# identical for both rows anyway.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = 19
$row2 = 20

$cols = @(
    "C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V",
    "W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM",
    "AN","AO","AP","AQ","AR","AS","AT","AU","AV","AW","AX","AY","AZ","BA","BB","BC",
    "BD","BE","BF","BG","BH","BI","BJ","BK","BL","BM","BN","BO","BP","BQ","BR","BS",
    "BT","BU","BV","BW","BX","BY","BZ","CA","CB","CC","CD","CE","CF","CG","CH","CI",
    "CJ","CK","CL","CM","CN","CO","CP","CQ","CR","CS","CT","CU","CV","CW","CX","CY",
    "CZ","DA","DB","DC","DD","DE","DF","DG","DH","DI","DJ"
)

foreach ($col in $cols) {
    $cell1 = $ws.Range($col + $row1)
    $cell2 = $ws.Range($col + $row2)
    $tmp = $cell1.Value2
    $cell1.Value2 = $cell2.Value2
    $cell2.Value2 = $tmp
}
